## Slide 19 ("Gracias" / contact info slide): turn the e-mail address into a
## mailto: hyperlink and append two new lines with office-hours info below it.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(19)
$sh = $s.Shapes.Item(2)              # "Marcador de texto 2" - name / e-mail box
$tr = $sh.TextFrame.TextRange

## 1) Append the two new paragraphs first, while the trailing run still has no
##    hyperlink, so the new runs pick up plain (non-linked) formatting that
##    matches the rest of the textbox.
$tr.InsertAfter("`rM `t– 16:00 – 18:00") | Out-Null
$tr.InsertAfter("`rJV `t- 14:00 -16:00") | Out-Null

## 2) Turn the e-mail address (the original second paragraph, still the first
##    23 characters after the implicit paragraph break at position 26) into a
##    mailto: hyperlink.
$emailRange = $tr.Characters(27, 23)
$emailRange.ActionSettings(1).Hyperlink.Address = "mailto:juan.salasf@udea.edu.co"
